$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8695881366729736
$ws.Range("B1").Value = 3.601375102996826
$ws.Range("C1").Value = 1.955960512161255
$ws.Range("D1").Value = 1.075894951820374
$ws.Range("E1").Value = 1.13348126411438
